$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct the contact data (per the final upload) ---
# Row 2 (first contact)
$ws.Range("B2").Value = "rohan"
$ws.Range("J2").Value = "tintu@gmail.com"
$ws.Range("L2").Value = "Reading ,Drawing"

# Row 3 (second contact)
$ws.Range("B3").Value = "mini"
$ws.Range("H3").Value = "abcd"
$ws.Range("L3").Value = "Reading ,Writing"

# --- Make the Pincode / Phone font color explicit black instead of the theme color ---
$ws.Range("I1:I3").Font.Color = 0
$ws.Range("K1:K3").Font.Color = 0

# --- Slightly taller header/data rows ---
$ws.Range("A1:A3").EntireRow.RowHeight = 19.5
